$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up a few header/summary labels: add a missing trailing period and
# capitalize the first letter.
$ws.Range("J2").Value = "Штраф, руб."          # was "Штраф, руб"

$ws.Range("B40").Value = "Общая сумма, руб."                      # was "общая сумма, руб."
$ws.Range("B41").Value = "Средняя площадь, кв.м."                 # was "средняя площадь, кв.м."
$ws.Range("B42").Value = "Максимальный срок просрочки, дней"      # was "максимальный срок просрочки, дней"
$ws.Range("B43").Value = "Максимальная сумма к оплате, руб."      # was "максимальная сумма к оплате, руб."

# Scroll the view back to the top-left (was parked at A16) and move the
# active selection to H2 (was B43).
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$ws.Range("H2").Select()
